$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Data")
$ws.Name = "Data table"
$ws.Activate()
